$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for the first server entry.
# Columns: A=ID, B=ServerID, C=Name, D=MaxOnline, E=CpuCount, F=IP, G=Port
# A2, B2, C2 and F2 need the "Text" number format (style index 1, numFmtId 49)
# so values like the leading-zero ServerID survive as text.
$ws.Range("A2:C2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"

$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "GameServer_1"
$ws.Range("C2").Value = "GameServer_1"
$ws.Range("B2").Value = "000104001"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 4001

# The list validation on column F should no longer cover the header-adjacent
# row now that F2 carries real data - it now starts at F3.
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F3:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Move the active selection to G3 (single cell), matching the saved view state.
$ws.Range("G3").Select() | Out-Null
